# Generate Report for Handoff
#
# The 45f196f5-c961-47f6-a64c-0fe8e7e03c1c item finished handoff later than
# e5ad083c/e8c8613a/01815df5, so its row moves from position 5 down to
# position 8 on every sheet (Overview, zh-cn, de-de), with the three rows in
# between shifting up by one. The "Latest Handoff" timestamp for
# 45f196f5-c961-47f6-a64c-0fe8e7e03c1c is also refreshed to reflect the later
# handoff time.

$wb = $excel.ActiveWorkbook

function Set-CellAndHyperlink {
    param(
        $ws,
        [string]$cellAddr,
        [string]$newValue,
        [bool]$isHyperlink
    )

    $ws.Range($cellAddr).Value = $newValue

    if ($isHyperlink) {
        foreach ($hl in $ws.Hyperlinks) {
            $addr = $hl.Range.Address()
            if ($addr -eq ('$' + $cellAddr.Replace(":", ":$"))) {
                $hl.TextToDisplay = $newValue
            }
        }
    }
}

# ---------------------------------------------------------------------
# Overview sheet: columns A (hyperlink), B, C, D
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

Set-CellAndHyperlink $wsOverview "A5" "e5ad083c-cc00-4e35-88ca-793ea5ae04d4.md" $true
$wsOverview.Range("B5").Value = "In Translation"
$wsOverview.Range("C5").Value = "In Translation"
$wsOverview.Range("D5").Value = "2016-39-13 12:39:21"

Set-CellAndHyperlink $wsOverview "A6" "e8c8613a-8c97-406c-bcc2-365d65fa4e12.md" $true
$wsOverview.Range("B6").Value = "In Translation"
$wsOverview.Range("C6").Value = "In Translation"
$wsOverview.Range("D6").Value = "2016-39-13 12:39:55"

Set-CellAndHyperlink $wsOverview "A7" "01815df5-c612-4267-8e3e-93304d033164.md" $true
$wsOverview.Range("B7").Value = "Ready for handoff"
$wsOverview.Range("C7").Value = "Ready for handoff"
$wsOverview.Range("D7").Value = "2016-40-13 12:40:13"

Set-CellAndHyperlink $wsOverview "A8" "45f196f5-c961-47f6-a64c-0fe8e7e03c1c.md" $true
$wsOverview.Range("B8").Value = "Ready for handoff"
$wsOverview.Range("C8").Value = "Ready for handoff"
$wsOverview.Range("D8").Value = "2016-42-13 12:42:46"

# ---------------------------------------------------------------------
# zh-cn sheet: columns A (hyperlink), C, D (hyperlink), E
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

Set-CellAndHyperlink $wsZhCn "A5" "e5ad083c-cc00-4e35-88ca-793ea5ae04d4.md" $true
$wsZhCn.Range("C5").Value = "In Translation"
Set-CellAndHyperlink $wsZhCn "D5" "e5ad083c-cc00-4e35-88ca-793ea5ae04d4.a1824740580a486af9d679d06fb3762408a8e21c.zh-cn.xlf" $true
$wsZhCn.Range("E5").Value = "2016-03-13 12:39:18"

Set-CellAndHyperlink $wsZhCn "A6" "e8c8613a-8c97-406c-bcc2-365d65fa4e12.md" $true
$wsZhCn.Range("C6").Value = "In Translation"
Set-CellAndHyperlink $wsZhCn "D6" "e8c8613a-8c97-406c-bcc2-365d65fa4e12.89390d9e1ef74186568c3e876084d2bb6b13a335.zh-cn.xlf" $true
$wsZhCn.Range("E6").Value = "2016-03-13 12:39:51"

Set-CellAndHyperlink $wsZhCn "A7" "01815df5-c612-4267-8e3e-93304d033164.md" $true
$wsZhCn.Range("C7").Value = "Ready for handoff"
Set-CellAndHyperlink $wsZhCn "D7" "01815df5-c612-4267-8e3e-93304d033164.d30c03fec02d87c3e414aaaa0841c0f031e041b6.zh-cn.xlf" $true
$wsZhCn.Range("E7").Value = "2016-03-13 12:40:09"

Set-CellAndHyperlink $wsZhCn "A8" "45f196f5-c961-47f6-a64c-0fe8e7e03c1c.md" $true
$wsZhCn.Range("C8").Value = "Ready for handoff"
Set-CellAndHyperlink $wsZhCn "D8" "45f196f5-c961-47f6-a64c-0fe8e7e03c1c.66ab83d27c5a3a936a3c8d6b71ce6a6c8c5f1d5d.zh-cn.xlf" $true
$wsZhCn.Range("E8").Value = "2016-03-13 12:42:43"

# ---------------------------------------------------------------------
# de-de sheet: columns A (hyperlink), C, D (hyperlink), E
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

Set-CellAndHyperlink $wsDeDe "A5" "e5ad083c-cc00-4e35-88ca-793ea5ae04d4.md" $true
$wsDeDe.Range("C5").Value = "In Translation"
Set-CellAndHyperlink $wsDeDe "D5" "e5ad083c-cc00-4e35-88ca-793ea5ae04d4.a1824740580a486af9d679d06fb3762408a8e21c.de-de.xlf" $true
$wsDeDe.Range("E5").Value = "2016-03-13 12:39:21"

Set-CellAndHyperlink $wsDeDe "A6" "e8c8613a-8c97-406c-bcc2-365d65fa4e12.md" $true
$wsDeDe.Range("C6").Value = "In Translation"
Set-CellAndHyperlink $wsDeDe "D6" "e8c8613a-8c97-406c-bcc2-365d65fa4e12.89390d9e1ef74186568c3e876084d2bb6b13a335.de-de.xlf" $true
$wsDeDe.Range("E6").Value = "2016-03-13 12:39:55"

Set-CellAndHyperlink $wsDeDe "A7" "01815df5-c612-4267-8e3e-93304d033164.md" $true
$wsDeDe.Range("C7").Value = "Ready for handoff"
Set-CellAndHyperlink $wsDeDe "D7" "01815df5-c612-4267-8e3e-93304d033164.d30c03fec02d87c3e414aaaa0841c0f031e041b6.de-de.xlf" $true
$wsDeDe.Range("E7").Value = "2016-03-13 12:40:13"

Set-CellAndHyperlink $wsDeDe "A8" "45f196f5-c961-47f6-a64c-0fe8e7e03c1c.md" $true
$wsDeDe.Range("C8").Value = "Ready for handoff"
Set-CellAndHyperlink $wsDeDe "D8" "45f196f5-c961-47f6-a64c-0fe8e7e03c1c.66ab83d27c5a3a936a3c8d6b71ce6a6c8c5f1d5d.de-de.xlf" $true
$wsDeDe.Range("E8").Value = "2016-03-13 12:42:46"
